# Regenerate the "K" column (strikeouts, column G) values in the save_data
# worksheet. The original values were derived from a different (incorrect)
# source; this updates them to the correct K counts, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 9
    4  = 5
    5  = 8
    6  = 10
    7  = 6
    8  = 2
    9  = 3
    10 = 9
    11 = 6
    12 = 8
    13 = 3
    14 = 3
    15 = 1
    16 = 6
    17 = 2
    18 = 9
    19 = 10
    20 = 3
    21 = 4
    22 = 2
    23 = 1
    24 = 4
    25 = 5
    26 = 1
    27 = 2
    28 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
